# Update the monitoring data: re-sort the empadronador totals and append
# two new rows (11 and 12) with additional records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2 through 12 (column A = empadronador, column B = total_registros)
$data = @(
    @("NIMA CARMEN KAREN DEL MILAGRO", 79),
    @("ARRUNATEGUI ESPINOZA JOVANNY", 77),
    @("ALZAMORA CHERRES SIRLEY YASMIN", 77),
    @("CARRION LAZARO MICHAEL LUIS", 76),
    @("PAZ ANASTACIO JUANITA ROSA", 76),
    @("ESPINOZA VALDIVIEZO JUNIOR RICARDO", 73),
    @("PULACHE LAZO VILMA YOHANA", 71),
    @("NAVARRO JUAREZ LIDIA", 70),
    @("LILIAN ROXANA VEGA GARCÍA", 63),
    @("DOMINGUEZ CUEVA MERLING DEL JESUS YOLINDA", 44),
    @("NIMA CRUZ ANA GRACIELA", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
